$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the wording of the "darker font / photos" system requirement in E5
# (System Requirement #1 for the "Website Pages" row) to use clearer,
# higher-contrast language.
$ws.Range("E5").Value = "Each site page shall use relevant articles typed in darker font with higher contrast and displayed next to their corresponding photos.  Currently, photos are placed on a page without corresponding label or reference.  Light font is used on a light gray background which can be difficult to read."
